$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "ff"
$ws.Range("F12").Value = "gg"
$ws.Range("C15").Value = "we"
$ws.Range("K16").Value = "nngf"

$ws.Range("K16").Select()
